$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A1").Value = "NameEn"
$ws.Range("B1").Value = "NameAr"

$ws.Range("B1").Select()
